$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'97.696.80"
$ws.Range("E2").Value = "  +1.29%  "
$ws.Range("D3").Value = "'3.721.84"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +13.31%  "
$ws.Range("D6").Value = "'238.61"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "'657.38"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E8").Value = "  +5.26%  "
$ws.Range("E9").Value = "  +5.91%  "
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").Value = "'3.720.43"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "'45.66"
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("D13").Value = "'0.0000313"
$ws.Range("E13").Value = "  +16.74%  "
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "'4.414.31"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "'97.446.31"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").Value = "'9.28"
$ws.Range("E18").Value = "  +3.29%  "
$ws.Range("D19").Value = "'3.727.68"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").Value = "'13.13"
$ws.Range("E20").Value = "  +2.82%  "
$ws.Range("D21").Value = "'18.97"
$ws.Range("E21").Value = "  -0.86%  "
$ws.Range("E22").Value = "  +2.69%  "
$ws.Range("D23").Value = "'533.04"
$ws.Range("E23").Value = "  +1.62%  "
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("E25").Value = "  +10.57%  "
$ws.Range("D26").Value = "'119.54"
$ws.Range("E26").Value = "  +16.50%  "
$ws.Range("D27").Value = "'6.93"
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("D28").Value = "'0.216"
$ws.Range("E28").Value = "  +28.17%  "
$ws.Range("D29").Value = "'13.48"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").Value = "'12.84"
$ws.Range("E30").Value = "  +2.90%  "
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").Value = "'0.193"
$ws.Range("E33").Value = "  +3.62%  "
$ws.Range("E34").Value = "  -3.25%  "
$ws.Range("D35").Value = "'33.23"
$ws.Range("E35").Value = "  +1.33%  "
$ws.Range("D36").Value = "'0.998"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("D37").Value = "'0.606"
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("D38").Value = "'641.52"
$ws.Range("E38").Value = "  -4.44%  "
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  +5.10%  "
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'41.46"
$ws.Range("E43").Value = "  +3.33%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "'0.494"
$ws.Range("E44").Value = "  +13.37%  "
$ws.Range("E45").Value = "  +2.73%  "
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("E48").Value = "  +3.95%  "
$ws.Range("D49").Value = "'8.98"
$ws.Range("D50").Value = "'23.67"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("E51").Value = "  +5.78%  "
